$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.131.58"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "1.666.59"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  -0.43%  "

$ws.Range("D5").Value = "'209.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.66%  "

$ws.Range("D6").Value = "'0.5204"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.13%  "

$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").Value = "'0.2615"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.50%  "

$ws.Range("D9").Value = "'0.06319"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.28%  "

$ws.Range("D10").Value = "'21.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.54%  "

$ws.Range("D11").Value = "'0.07525"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.03%  "

$ws.Range("D12").Value = "1.669.43"
$ws.Range("E12").Value = "  -1.14%  "

$ws.Range("D13").Value = "'4.423"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.24%  "

$ws.Range("D14").Value = "'0.5479"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.25%  "

$ws.Range("D15").Value = "'66.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("D16").Value = "'0.000007936"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.04%  "

$ws.Range("D17").Value = "26.170.32"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D19").Value = "'4.723"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.38%  "

$ws.Range("D20").Value = "'186.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.79%  "

$ws.Range("D21").Value = "'10.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.34%  "

$ws.Range("D22").Value = "'6.172"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.27%  "

$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("D24").Value = "'148.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("E25").Value = "  -1.95%  "

$ws.Range("D26").Value = "'7.477"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.48%  "

$ws.Range("D27").Value = "'15.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").Value = "'0.06353"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("D30").Value = "'1.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.74%  "

$ws.Range("D31").Value = "'3.495"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.83%  "

$ws.Range("D32").Value = "'3.407"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.66%  "

$ws.Range("D33").Value = "'1.640"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.87%  "

$ws.Range("E34").Value = "  -2.38%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.407"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.00%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6006"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.73%  "

$ws.Range("E37").Value = "  -0.25%  "

$ws.Range("D38").Value = "1.109.17"
$ws.Range("E38").Value = "  +0.17%  "

$ws.Range("D39").Value = "'6.097"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.22%  "

$ws.Range("D40").Value = "'0.01613"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("D41").Value = "'0.8654"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.01%  "

$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("D43").Value = "'100.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.90%  "

$ws.Range("D44").Value = "1.819.47"
$ws.Range("E44").Value = "  -1.00%  "

$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.82%  "

$ws.Range("D46").Value = "'55.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.16%  "

$ws.Range("D47").Value = "'0.9987"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.89%  "

$ws.Range("D48").Value = "'8.043"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").Value = "'0.05230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("D50").Value = "'0.4246"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("D51").Value = "'5.914"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.21%  "
